$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting the existing rows 20-41 down to 21-42
# (keeping all of their original values intact).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44987
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100101
$ws.Range("H20").Value = "Berries"
$ws.Range("I20").Value = 100101006
$ws.Range("J20").Value = "Higo"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = 18000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 18000
$ws.Range("Q20").Value = "$/bandeja 7 kilos"
$ws.Range("R20").Value = "Provincia de Santiago"
$ws.Range("S20").Value = 2571
$ws.Range("T20").Value = 7
